$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$note = "have decreased sizes of part at first census, but this has no effect. Wondering what I'm doing wrong."

# Rows that get a brand-new J-column note (yellow highlighted, like the existing notes column)
$newNoteRows = @(2,3,4,5,6,8,9)
foreach ($r in $newNoteRows) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Value = $note
    $cell.Interior.Color = 65535
}

# Rows whose existing J-column note just gets the yellow highlight applied (value unchanged)
$highlightOnlyRows = @(11,12,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
foreach ($r in $highlightOnlyRows) {
    $ws.Cells.Item($r, 10).Interior.Color = 65535
}

# Update the active selection shown in the sheet view
$ws.Range("F16").Select() | Out-Null
